$d = $word.ActiveDocument

# 1. "...sporting all capabilities but not very great numbers..."
#    -> "...sporting all capabilities but not in very great numbers..."
$null = $d.Content.Find.Execute(
    "sporting all capabilities but not very great numbers",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "sporting all capabilities but not in very great numbers", 2)

# 2. "intercepting contraband" -> "interdicting contraband"
$null = $d.Content.Find.Execute(
    "intercepting contraband",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "interdicting contraband", 2)

# 3. Remove the stray _GoBack bookmark left over from editing.
if ($d.Bookmarks.Exists("_GoBack")) {
    $null = $d.Bookmarks("_GoBack").Delete()
}
